# Weekly fruit/vegetable price update: insert a new daily record for
# "Vega Monumental Concepción" / Mango as row 93, pushing the previously
# existing rows 93:121 down to 94:122.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 93 (existing rows 93-121 shift down to 94-122).
$ws.Rows("93:93").Insert()

# Populate the newly inserted row with the new weekly price entry.
$ws.Cells.Item(93, 1).Value = 11
$ws.Cells.Item(93, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(93, 3).Value = "Bíobío"
$ws.Cells.Item(93, 4).Value = 44782
$ws.Cells.Item(93, 5).Value = 8
$ws.Cells.Item(93, 6).Value = "Fruta"
$ws.Cells.Item(93, 7).Value = 100108
$ws.Cells.Item(93, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(93, 9).Value = 100108002
$ws.Cells.Item(93, 10).Value = "Mango"
$ws.Cells.Item(93, 11).Value = "Sin especificar"
$ws.Cells.Item(93, 12).Value = "Primera"
$ws.Cells.Item(93, 13).Value = 190
$ws.Cells.Item(93, 14).Value = 9000
$ws.Cells.Item(93, 15).Value = 10000
$ws.Cells.Item(93, 16).Value = 9526
$ws.Cells.Item(93, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(93, 18).Value = "Perú"
$ws.Cells.Item(93, 19).Value = 2382
$ws.Cells.Item(93, 20).Value = 4
